$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To do")

# Insert a new row above row 23 (pushes existing rows 23.. down by one)
$ws.Rows.Item(23).EntireRow.Insert()

# Fill in the new row 23 with the new to-do item
$ws.Cells.Item(23, 1).Value = "NETWORK IMPORT"
$ws.Cells.Item(23, 2).Value = "Check error handling when importing multi-arm trials from pairs format if some contrasts are missing"
$ws.Cells.Item(23, 3).Value = "me"
$ws.Cells.Item(23, 4).Value = 43356
$ws.Cells.Item(23, 4).NumberFormat = "m/d/yy"

# Match source row styling (row above uses the same formats)
$ws.Rows.Item(22).Copy()
$ws.Rows.Item(23).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-apply the values since PasteSpecial(formats) shouldn't touch them, but ensure correctness
$ws.Cells.Item(23, 1).Value = "NETWORK IMPORT"
$ws.Cells.Item(23, 2).Value = "Check error handling when importing multi-arm trials from pairs format if some contrasts are missing"
$ws.Cells.Item(23, 3).Value = "me"
$ws.Cells.Item(23, 4).Value = 43356

$ws.Rows.Item(23).RowHeight = 30
